# TC09 Canine Study workbook - "updated 2 icdc scripts to resolve wait time issue"
#
# The "CasesTab" row's Cypher query (cell B2 on the "startup" sheet) is
# rewritten to drop the trailing `cohort` lookup (the
# "OPTIONAL MATCH (co:cohort)..." *return* column), which is what was
# causing the slow query / wait-time issue. The query still declares
# `co` in its WITH clause (unused cohort param kept out of RETURN).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newCaseQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['UBC01'] and demo.breed in ['Mixed Breed', 'Scottish Terrier','Shetland Sheepdog']and diag.disease_term in ['Bladder Cancer','Healthy Control'] and diag.primary_disease_site in ['Bladder']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $newCaseQuery

# Row 2 wraps B2/C2 (style index 1 => wrapText). With the shorter query
# text the row no longer needs to be as tall as before (304.5 -> 290,
# same height the SamplesTab row below already uses).
$ws.Rows(2).RowHeight = 290

# Reflect where the author was working: scrolled back to the top of the
# sheet with the just-edited CasesTab query cell selected.
$ws.Range("B2").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
